$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update / extend the employee list (Name, Image URL) in columns A:B.
# Values are written in the same order the original author typed them in
# Excel so the resulting shared-strings table matches exactly (a couple of
# rows were filled out of the usual A-then-B order).

$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Image URL"

$ws.Cells.Item(2, 1).Value = "abhishek"
$ws.Cells.Item(2, 2).Value = "meme_images/abhishek.png"

$ws.Cells.Item(3, 1).Value = "azharudhin"
$ws.Cells.Item(3, 2).Value = "meme_images/azhar.png"

$ws.Cells.Item(4, 1).Value = "archana"
$ws.Cells.Item(4, 2).Value = "meme_images/archana.png"

$ws.Cells.Item(5, 1).Value = "akash"
$ws.Cells.Item(5, 2).Value = "meme_images/akash_deep.png"

$ws.Cells.Item(6, 1).Value = "adrija"
$ws.Cells.Item(6, 2).Value = "meme_images/adrija.png"

$ws.Cells.Item(7, 1).Value = "bishal"
$ws.Cells.Item(7, 2).Value = "meme_images/bishal.png"

$ws.Cells.Item(8, 1).Value = "daksh"
$ws.Cells.Item(8, 2).Value = "meme_images/Daksh.png"

$ws.Cells.Item(9, 1).Value = "deepthi nair"
$ws.Cells.Item(9, 2).Value = "meme_images/deepthi_nair.png"

$ws.Cells.Item(10, 1).Value = "gishika"
$ws.Cells.Item(10, 2).Value = "meme_images/gishika.png"

$ws.Cells.Item(11, 1).Value = "harsh"
$ws.Cells.Item(12, 1).Value = "deepthi valsan"
$ws.Cells.Item(12, 2).Value = "meme_images/deepthi_valsan.png"
$ws.Cells.Item(11, 2).Value = "meme_images/harsh.png"

$ws.Cells.Item(13, 1).Value = "zain"
$ws.Cells.Item(13, 2).Value = "meme_images/zain.png"

$ws.Cells.Item(14, 1).Value = "vijay"
$ws.Cells.Item(14, 2).Value = "meme_images/vijay.png"

$ws.Cells.Item(15, 1).Value = "sneha"
$ws.Cells.Item(15, 2).Value = "meme_images/sneha.png"

$ws.Cells.Item(16, 1).Value = "shamim"
$ws.Cells.Item(16, 2).Value = "meme_images/shamim.png"

$ws.Cells.Item(17, 1).Value = "sarwesh"
$ws.Cells.Item(17, 2).Value = "meme_images/sarwesh.png"

$ws.Cells.Item(18, 1).Value = "sanjana"
$ws.Cells.Item(18, 2).Value = "meme_images/sanjana.png"

$ws.Cells.Item(19, 1).Value = "ragav"
$ws.Cells.Item(19, 2).Value = "meme_images/ragav.png"

$ws.Cells.Item(20, 1).Value = "prashant_singhal"
$ws.Cells.Item(20, 2).Value = "meme_images/prashant_singhal.png"

$ws.Cells.Item(21, 1).Value = "pavani"
$ws.Cells.Item(21, 2).Value = "meme_images/pavani.png"

$ws.Cells.Item(22, 1).Value = "parashant"
$ws.Cells.Item(22, 2).Value = "meme_images/parashant.png"

$ws.Cells.Item(23, 1).Value = "lipika"
$ws.Cells.Item(23, 2).Value = "meme_images/lipika.png"

$ws.Cells.Item(24, 2).Value = "meme_images/imad.png"
$ws.Cells.Item(24, 1).Value = "imad"

# Auto-fit column A so it best-fits the new, longer names (e.g. "deepthi valsan")
$ws.Columns.Item(1).AutoFit() | Out-Null

# Scroll / selection state as captured in the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("G21").Select()
